# Rename the single worksheet from "Property1" to "DataNode" to unify the
# conception of DataNode / DataTable / Entity across the data config sheets.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# The workbook was last saved with the selection sitting on D26 (the sheet
# view's frozen pane keeps the header rows visible while the cursor was
# moved further down/right before the file was saved).
$ws.Range("D26").Select()
